$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 28 (2023Q3) with the refreshed BIBI metrics
$ws.Range("C28").Value = 423
$ws.Range("D28").Value = 48
$ws.Range("E28").Value = 375
$ws.Range("F28").Value = 7.476635514018691
